# Generate Report for Handback
# Updates timestamps / status text produced by a new handback report run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for 5eb8ba03...md / d7a7430f...md rows
$wsOverview.Range("G4").Value = "2016-08-31 06:18:20"
$wsOverview.Range("G5").Value = "2016-08-31 06:18:20"

# zh-cn sheet: Priority changed from "ht" to "mt" for the 5eb8ba03 / d7a7430f rows
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"

# zh-cn sheet: Correspond Handoff Datetime
$wsZhCn.Range("H4").Value = "2016-08-31 06:18:15"
$wsZhCn.Range("H5").Value = "2016-08-31 06:18:15"

# zh-cn sheet: Correspond Handback DateTime
$wsZhCn.Range("K4").Value = "2016-08-31 06:18:32"
$wsZhCn.Range("K5").Value = "2016-08-31 06:18:32"

# de-de sheet: Priority changed from "ht" to "mt" for the 5eb8ba03 / d7a7430f rows
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"

# de-de sheet: Correspond Handoff Datetime (shares the same timestamp as Overview G4/G5)
$wsDeDe.Range("H4").Value = "2016-08-31 06:18:20"
$wsDeDe.Range("H5").Value = "2016-08-31 06:18:20"

# de-de sheet: Correspond Handback DateTime
$wsDeDe.Range("K4").Value = "2016-08-31 06:18:39"
$wsDeDe.Range("K5").Value = "2016-08-31 06:18:39"
